$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark.  In the original document
#    it sits at the very start of the "Currently when you enable/disable..."
#    list item (paragraph 3); it needs to move to the end of the new
#    bullet point we are about to add.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Add a brand-new bullet point right after the "Currently when you
#    enable/disable..." paragraph, using the same list formatting
#    (ListParagraph style, same numbering).
# ------------------------------------------------------------------
$sourcePara = $d.Paragraphs(3)
$sourcePara.Range.InsertParagraphAfter()

$firstPara = $d.Paragraphs(4)
$firstInsertionPoint = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)
$firstInsertionPoint.InsertAfter("Need to have the bullet points on the error window wrap to a new line inside the error window if the screen becomes two small.")

# Add the second sentence as its own run by typing it into a temporary
# paragraph directly after the first one, then merging the two
# paragraphs back together (deleting the paragraph mark keeps the two
# pieces of text as separate <w:r> runs instead of collapsing them into
# a single run).
$d.Paragraphs(4).Range.InsertParagraphAfter()
$secondPara = $d.Paragraphs(5)
$secondInsertionPoint = $d.Range($secondPara.Range.Start, $secondPara.Range.Start)
$secondInsertionPoint.InsertAfter("  This might be because I am not having it stretch to fill the stack panel, but I am unsure about this.")

$paraMarkRange = $d.Range($d.Paragraphs(4).Range.End - 1, $d.Paragraphs(4).Range.End)
$paraMarkRange.Delete()

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark, collapsed, at the end of the
#    new bullet point (after the text, before the paragraph mark).
#    A collapsed range exactly at (paragraph.End - 1) confuses bookmark
#    placement, so a temporary character is inserted first to give the
#    bookmark a stable, non-edge position; the character is removed
#    again once the bookmark has been created.
# ------------------------------------------------------------------
$bulletPara = $d.Paragraphs(4)
$endPos = $bulletPara.Range.End - 1
$tempRange = $d.Range($endPos, $endPos)
$tempRange.InsertAfter("X")

$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$tempCharRange = $d.Range($endPos, $endPos + 1)
$tempCharRange.Text = ""
